$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("G2").Value = 2.374260666666666
$ws.Range("H2").Value = 7.122781999999999
$ws.Range("I2").Value = 0.0276017086472712
$ws.Range("J2").Value = 0.0276017086472712
$ws.Range("M2").Value = 52.91852733333334
$ws.Range("N2").Value = 158.755582
$ws.Range("O2").Value = 0.9912603569328422
$ws.Range("P2").Value = 0.9912603569328421
$ws.Range("Q2").Value = 125.6423779854582
$ws.Range("R2").Value = 1130.781401869124
$ws.Range("S2").Value = 0.02736047956565037
$ws.Range("T2").Value = 0.02736047956565036
$ws.Range("G3").Value = 2.374260666666666
$ws.Range("H3").Value = 7.122781999999999
$ws.Range("I3").Value = 0.0276017086472712
$ws.Range("J3").Value = 0.0276017086472712
$ws.Range("O3").Value = 0.003851187374513192
$ws.Range("P3").Value = 0.003851187374513192
$ws.Range("Q3").Value = 0.4881384960239999
$ws.Range("R3").Value = 4.393246464215999
$ws.Range("S3").Value = 0.0001062993518573625
$ws.Range("T3").Value = 0.0001062993518573624
$ws.Range("G4").Value = 2.374260666666666
$ws.Range("H4").Value = 7.122781999999999
$ws.Range("I4").Value = 0.0276017086472712
$ws.Range("J4").Value = 0.0276017086472712
$ws.Range("K4").Value = 3
$ws.Range("L4").Value = 1
$ws.Range("M4").Value = 0.2609706666666667
$ws.Range("N4").Value = 0.7829120000000001
$ws.Range("O4").Value = 0.004888455692644593
$ws.Range("P4").Value = 0.004888455692644592
$ws.Range("Q4").Value = 0.6196123890204444
$ws.Range("R4").Value = 5.576511501184
$ws.Range("S4").Value = 0.0001349297297634704
$ws.Range("T4").Value = 0.0001349297297634704
$ws.Range("I5").Value = 0.8942818522422411
$ws.Range("J5").Value = 0.8942818522422411
$ws.Range("M5").Value = 52.91852733333334
$ws.Range("N5").Value = 158.755582
$ws.Range("O5").Value = 0.9912603569328422
$ws.Range("P5").Value = 0.9912603569328421
$ws.Range("Q5").Value = 4070.751558927987
$ws.Range("R5").Value = 36636.76403035189
$ws.Range("S5").Value = 0.8864661480522071
$ws.Range("T5").Value = 0.886466148052207
$ws.Range("I6").Value = 0.8942818522422411
$ws.Range("J6").Value = 0.8942818522422411
$ws.Range("O6").Value = 0.003851187374513192
$ws.Range("P6").Value = 0.003851187374513192
$ws.Range("S6").Value = 0.003444046978611591
$ws.Range("T6").Value = 0.003444046978611591
$ws.Range("I7").Value = 0.8942818522422411
$ws.Range("J7").Value = 0.8942818522422411
$ws.Range("K7").Value = 3
$ws.Range("L7").Value = 1
$ws.Range("M7").Value = 0.2609706666666667
$ws.Range("N7").Value = 0.7829120000000001
$ws.Range("O7").Value = 0.004888455692644593
$ws.Range("P7").Value = 0.004888455692644592
$ws.Range("Q7").Value = 20.07513817374578
$ws.Range("R7").Value = 180.676243563712
$ws.Range("S7").Value = 0.004371657211422334
$ws.Range("T7").Value = 0.004371657211422333
$ws.Range("G8").Value = 6.625048
$ws.Range("H8").Value = 19.875144
$ws.Range("I8").Value = 0.07701877356495823
$ws.Range("J8").Value = 0.07701877356495825
$ws.Range("M8").Value = 52.91852733333334
$ws.Range("N8").Value = 158.755582
$ws.Range("O8").Value = 0.9912603569328422
$ws.Range("P8").Value = 0.9912603569328421
$ws.Range("Q8").Value = 350.5877836726453
$ws.Range("R8").Value = 3155.290053053808
$ws.Range("S8").Value = 0.07634565697453025
$ws.Range("T8").Value = 0.07634565697453026
$ws.Range("G9").Value = 6.625048
$ws.Range("H9").Value = 19.875144
$ws.Range("I9").Value = 0.07701877356495823
$ws.Range("J9").Value = 0.07701877356495825
$ws.Range("O9").Value = 0.003851187374513192
$ws.Range("P9").Value = 0.003851187374513192
$ws.Range("Q9").Value = 1.362083368608
$ws.Range("R9").Value = 12.258750317472
$ws.Range("S9").Value = 0.0002966137283538576
$ws.Range("T9").Value = 0.0002966137283538576
$ws.Range("G10").Value = 6.625048
$ws.Range("H10").Value = 19.875144
$ws.Range("I10").Value = 0.07701877356495823
$ws.Range("J10").Value = 0.07701877356495825
$ws.Range("K10").Value = 3
$ws.Range("L10").Value = 1
$ws.Range("M10").Value = 0.2609706666666667
$ws.Range("N10").Value = 0.7829120000000001
$ws.Range("O10").Value = 0.004888455692644593
$ws.Range("P10").Value = 0.004888455692644592
$ws.Range("Q10").Value = 1.728943193258667
$ws.Range("R10").Value = 15.560488739328
$ws.Range("S10").Value = 0.0003765028620741249
$ws.Range("T10").Value = 0.0003765028620741249
$ws.Range("G11").Value = 0.09441966666666667
$ws.Range("H11").Value = 0.283259
$ws.Range("I11").Value = 0.001097665545529457
$ws.Range("J11").Value = 0.001097665545529457
$ws.Range("M11").Value = 52.91852733333334
$ws.Range("N11").Value = 158.755582
$ws.Range("O11").Value = 0.9912603569328422
$ws.Range("P11").Value = 0.9912603569328421
$ws.Range("Q11").Value = 4.996549711304223
$ws.Range("R11").Value = 44.968947401738
$ws.Range("S11").Value = 0.001088072340454412
$ws.Range("T11").Value = 0.001088072340454412
$ws.Range("G12").Value = 0.09441966666666667
$ws.Range("H12").Value = 0.283259
$ws.Range("I12").Value = 0.001097665545529457
$ws.Range("J12").Value = 0.001097665545529457
$ws.Range("O12").Value = 0.003851187374513192
$ws.Range("P12").Value = 0.003851187374513192
$ws.Range("Q12").Value = 0.019412305788
$ws.Range("R12").Value = 0.174710752092
$ws.Range("S12").Value = 0.000004227315690381178
$ws.Range("T12").Value = 0.000004227315690381178
$ws.Range("G13").Value = 0.09441966666666667
$ws.Range("H13").Value = 0.283259
$ws.Range("I13").Value = 0.001097665545529457
$ws.Range("J13").Value = 0.001097665545529457
$ws.Range("K13").Value = 3
$ws.Range("L13").Value = 1
$ws.Range("M13").Value = 0.2609706666666667
$ws.Range("N13").Value = 0.7829120000000001
$ws.Range("O13").Value = 0.004888455692644593
$ws.Range("P13").Value = 0.004888455692644592
$ws.Range("Q13").Value = 0.02464076335644445
$ws.Range("R13").Value = 0.221766870208
$ws.Range("S13").Value = 0.000005365889384663304
$ws.Range("T13").Value = 0.000005365889384663303
